# "PLO Stats-this session" now reflects a different (newer) session: the
# player roster/stats for rows 2-5 change, the two trailing rows (Jacob,
# and the old duplicate-looking Scott row) disappear entirely, and the
# seven charts that plot this sheet need their category/value series
# ranges shrunk from row 7 down to row 5 to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("PLO Stats-this session")

# --- Row 2 (Raymond) : update stats in place ---------------------------
$ws.Range("B2").Value = 60
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = -60
$ws.Range("F2").Value = 0.713
$ws.Range("G2").Value = 0.438
$ws.Range("H2").Value = 0.05
$ws.Range("I2").Value = 0.497
$ws.Range("J2").Value = 0.175
$ws.Range("K2").Value = 0.075
$ws.Range("L2").Value = 2.31
$ws.Range("M2").Value = 18
$ws.Range("N2").Value = 40
$ws.Range("O2").Value = 72.64
$ws.Range("P2").Value = 59.64
$ws.Range("Q2").Value = 80
$ws.Range("R2").Value = 0.429

# --- Row 3 : was Fish, now Cedric's stats -------------------------------
$ws.Range("A3").Value = "Cedric"
$ws.Range("B3").Value = 50
$ws.Range("D3").Value = -50
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.614
$ws.Range("G3").Value = 0.014
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0.268
$ws.Range("J3").Value = 0.2
$ws.Range("K3").Value = 0.086
$ws.Range("L3").Value = 0.6
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 29.45
$ws.Range("P3").Value = 36.72
$ws.Range("Q3").Value = 70
$ws.Range("R3").Value = 0.429

# T3 holds a literal text date like "04/29/21" (not a real Excel date).
# Force Text format first so the new text isn't auto-converted into a
# serial date number / date-formatted cell, then clear the format change
# back off so no stray style sticks to the cell.
$ws.Range("T3").NumberFormat = "@"
$ws.Range("T3").Value = "06/10/21"
$ws.Range("T3").ClearFormats()

# --- Row 4 : was Cedric, now Fish's stats -------------------------------
$ws.Range("A4").Value = "Fish"
$ws.Range("B4").Value = 40
$ws.Range("C4").Value = 80.84
$ws.Range("D4").Value = 40.84
$ws.Range("F4").Value = 0.772
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.171
$ws.Range("J4").Value = 0.207
$ws.Range("K4").Value = 0.098
$ws.Range("L4").Value = 0.29
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 175.94
$ws.Range("P4").Value = 51.3
$ws.Range("Q4").Value = 92
$ws.Range("R4").Value = 0.474

# --- Row 5 : was Kynan, now Scott's stats -------------------------------
$ws.Range("A5").Value = "Scott"
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 89.16
$ws.Range("D5").Value = 69.16
$ws.Range("F5").Value = 0.788
$ws.Range("G5").Value = 0.288
$ws.Range("H5").Value = 0.0192
$ws.Range("I5").Value = 0.45
$ws.Range("J5").Value = 0.173
$ws.Range("K5").Value = 0.115
$ws.Range("L5").Value = 1.5
$ws.Range("M5").Value = 6
$ws.Range("N5").Value = 14
$ws.Range("O5").Value = 117.77
$ws.Range("P5").Value = 37.59
$ws.Range("Q5").Value = 52
$ws.Range("R5").Value = 0.667

# --- Rows 6-7 (Jacob, old Scott) no longer exist in this session -------
$ws.Range("A6:A7").EntireRow.Delete()

# --- Fix up the 7 charts on this sheet: $A$2:$A$7 -> $A$2:$A$5, and the
#     matching value-series column likewise $x$2:$x$7 -> $x$2:$x$5.
function Set-PloSeries($chartIndex, $seriesIndex, $headerCol, $valueCol) {
    $chart = $ws.ChartObjects().Item($chartIndex).Chart
    $sc = $chart.SeriesCollection()
    $s = $sc.Item($seriesIndex)
    $s.Formula = "=SERIES('PLO Stats-this session'!" + $headerCol + "1,'PLO Stats-this session'!`$A`$2:`$A`$5,'PLO Stats-this session'!`$" + $valueCol + "`$2:`$" + $valueCol + "`$5," + $seriesIndex + ")"
}

# Chart 1 (VPIP, Pre-flop raise, 3-bet): F, G, H
Set-PloSeries 1 1 "F" "F"
Set-PloSeries 1 2 "G" "G"
Set-PloSeries 1 3 "H" "H"

# Chart 2 (Aggression factor): L
Set-PloSeries 2 1 "L" "L"

# Chart 3 (C-bets / opportunities): M, N
Set-PloSeries 3 1 "M" "M"
Set-PloSeries 3 2 "N" "N"

# Chart 4 (Went to / won at showdown): J, K
Set-PloSeries 4 1 "J" "J"
Set-PloSeries 4 2 "K" "K"

# Chart 5 (WTSD rel): R
Set-PloSeries 5 1 "R" "R"

# Chart 6 (At / before showdown): O, P
Set-PloSeries 6 1 "O" "O"
Set-PloSeries 6 2 "P" "P"

# Chart 7 (Hands played): Q
Set-PloSeries 7 1 "Q" "Q"
